$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refactor: rename the "id" column header to "Name"
$ws.Range("A1").Value = "Name"
